# Archetypes and Aspects - update shared strings and cell layout
# "Electric effects, and Thunder Cherries are on cherry trees now!"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename aspect "Glorious" -> "Commanding" (column K) ---
$ws.Range("K1").Value = "Commanding"

# --- Row 4: move the "ALL" marker from column H to column K ---
$ws.Range("H4").Value = $null
$ws.Range("K4").Value = "ALL"

# --- Row 8: rename archetype "Knight" -> "Dreadnought", and restructure
#     its aspect markers: remove "Royal Knight"/"White Knight", introduce
#     "Shogun" and "Mamluk" ---
$ws.Range("A8").Value = "Dreadnought"
$ws.Range("I8").Value = $null
$ws.Range("K8").Value = "Shogun"
$ws.Range("L8").Value = "Death Knight"
$ws.Range("M8").Value = "ALL"
$ws.Range("V8").Value = "Mamluk"

# --- New row 12: Mentalist archetype ---
$ws.Range("A12").Value = "Mentalist"
$ws.Range("B12").Value = "Ardent"
$ws.Range("C12").Value = "Deadeye"
$ws.Range("E12").Value = "Empath"
$ws.Range("K12").Value = "ALL"
$ws.Range("P12").Value = "ALL"
$ws.Range("Q12").Value = "ALL"

# --- New row 13: Exemplar archetype ---
$ws.Range("A13").Value = "Exemplar"
$ws.Range("D13").Value = "ALL"
$ws.Range("H13").Value = "ALL"
$ws.Range("J13").Value = "ALL"
$ws.Range("N13").Value = "Battle Savant"
$ws.Range("P13").Value = "Erased Spy"
$ws.Range("V13").Value = "Atom Lord"

# --- Update selection to follow the newly added rows ---
$ws.Range("A14").Select()
